$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 770.8125
$ws.Range("I33").Value = 617.61536
$ws.Range("J33").Value = 1434.6666
$ws.Range("K33").Value = 617.61536
$ws.Range("L33").Value = 1434.6666
$ws.Range("M33").Value = -388.61536
$ws.Range("N33").Value = -1892.6666
$ws.Range("H40").Value = 2560.25
$ws.Range("I40").Value = 2480.4
$ws.Range("J40").Value = 2640.1
$ws.Range("K40").Value = 2480.4
$ws.Range("L40").Value = 2640.1
$ws.Range("M40").Value = -2305.4
$ws.Range("N40").Value = -2990.1
$ws.Range("H51").Value = 5555
$ws.Range("I51").Value = 1332
$ws.Range("J51").Value = 7666.5
$ws.Range("K51").Value = 1332
$ws.Range("L51").Value = 7666.5
$ws.Range("M51").Value = -848
$ws.Range("N51").Value = -8634.5
$ws.Range("H55").Value = 280.57144
$ws.Range("I55").Value = 388.33334
$ws.Range("J55").Value = 199.75
$ws.Range("K55").Value = 388.33334
$ws.Range("L55").Value = 199.75
$ws.Range("M55").Value = -174.33334
$ws.Range("N55").Value = -627.75
$ws.Range("H62").Value = 32698.945
$ws.Range("I62").Value = 4149.1816
$ws.Range("J62").Value = 77562.86
$ws.Range("K62").Value = 4149.1816
$ws.Range("L62").Value = 77562.86
$ws.Range("M62").Value = -3525.1816
$ws.Range("N62").Value = -78810.86
$ws.Range("H65").Value = 32698.945
$ws.Range("I65").Value = 4149.1816
$ws.Range("J65").Value = 77562.86
$ws.Range("K65").Value = 20745.908
$ws.Range("L65").Value = 387814.3
$ws.Range("M65").Value = -17625.908
$ws.Range("N65").Value = -394054.3
$ws.Range("H80").Value = 52499
$ws.Range("I80").Value = 52499
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 157497
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -156499
$ws.Range("H83").Value = 52499
$ws.Range("I83").Value = 52499
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 472491
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -467499
$ws.Range("H116").Value = 25414198
$ws.Range("I116").Value = 15697522
$ws.Range("J116").Value = 47623744
$ws.Range("K116").Value = 15697522
$ws.Range("L116").Value = 47623744
$ws.Range("M116").Value = -15694080
$ws.Range("N116").Value = -47630628

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1186.25
$ws.Range("I2").Value = 1021.36365
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 1021.36365
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = -908.36365
$ws.Range("N2").Value = -3226
$ws.Range("H38").Value = 5616.7144
$ws.Range("I38").Value = 2329.5
$ws.Range("J38").Value = 9999.666999999999
$ws.Range("K38").Value = 2329.5
$ws.Range("L38").Value = 9999.666999999999
$ws.Range("M38").Value = -1862.5
$ws.Range("N38").Value = -10933.667
$ws.Range("H62").Value = 19999
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 19999
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").Value = 19999
$ws.Range("N62").Value = -21247
$ws.Range("H65").Value = 19999
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 19999
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").Value = 59997
$ws.Range("N65").Value = -66237
$ws.Range("H116").Value = 1186.25
$ws.Range("I116").Value = 1021.36365
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 1021.36365
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = 1272.63635
$ws.Range("N116").Value = -7588
$ws.Range("H124").Value = 64832.668
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 64832.668
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 64832.668
$ws.Range("N124").Value = -74652.66800000001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1186.25
$ws.Range("I3").Value = 1021.36365
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 1021.36365
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = -907.36365
$ws.Range("N3").Value = -3228
$ws.Range("H99").Value = 1833.0769
$ws.Range("I99").Value = 1569.1666
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 1569.1666
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -71.16660000000002
$ws.Range("N99").Value = -7996
$ws.Range("H105").Value = 1392.2195
$ws.Range("I105").Value = 1194.6333
$ws.Range("J105").Value = 1931.091
$ws.Range("K105").Value = 1194.6333
$ws.Range("L105").Value = 1931.091
$ws.Range("M105").Value = 552.3667
$ws.Range("N105").Value = -5425.091

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3147.6
$ws.Range("I31").Value = 1126.7
$ws.Range("J31").Value = 3515.0364
$ws.Range("K31").Value = 1126.7
$ws.Range("L31").Value = 3515.0364
$ws.Range("M31").Value = -831.7
$ws.Range("N31").Value = -4105.0364
$ws.Range("H34").Value = 3147.6
$ws.Range("I34").Value = 1126.7
$ws.Range("J34").Value = 3515.0364
$ws.Range("K34").Value = 1126.7
$ws.Range("L34").Value = 3515.0364
$ws.Range("M34").Value = -924.7
$ws.Range("N34").Value = -3919.0364

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2553
$ws.Range("I5").Value = 1919.7
$ws.Range("J5").Value = 4664
$ws.Range("K5").Value = 5759.1
$ws.Range("L5").Value = 13992
$ws.Range("M5").Value = -5647.1
$ws.Range("N5").Value = -14216
$ws.Range("H26").Value = 78.2
$ws.Range("I26").Value = 78.2
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 234.6
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 53.39999999999998
$ws.Range("H42").Value = 1000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 1000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 3000
$ws.Range("N42").Value = -4068
$ws.Range("H55").Value = 2634.4375
$ws.Range("I55").Value = 1191.8334
$ws.Range("J55").Value = 3500
$ws.Range("K55").Value = 3575.5002
$ws.Range("L55").Value = 10500
$ws.Range("M55").Value = -3398.5002
$ws.Range("N55").Value = -10854
$ws.Range("H62").Value = 5937.6924
$ws.Range("I62").Value = 2575.111
$ws.Range("J62").Value = 13503.5
$ws.Range("K62").Value = 7725.333
$ws.Range("L62").Value = 40510.5
$ws.Range("M62").Value = -7039.333
$ws.Range("N62").Value = -41882.5
$ws.Range("H63").Value = 16989.55
$ws.Range("I63").Value = 18576
$ws.Range("J63").Value = 7999.6665
$ws.Range("K63").Value = 55728
$ws.Range("L63").Value = 23998.9995
$ws.Range("M63").Value = -54979
$ws.Range("N63").Value = -25496.9995
$ws.Range("H65").Value = 5937.6924
$ws.Range("I65").Value = 2575.111
$ws.Range("J65").Value = 13503.5
$ws.Range("K65").Value = 23175.999
$ws.Range("L65").Value = 121531.5
$ws.Range("M65").Value = -19743.999
$ws.Range("N65").Value = -128395.5
$ws.Range("H66").Value = 16989.55
$ws.Range("I66").Value = 18576
$ws.Range("J66").Value = 7999.6665
$ws.Range("K66").Value = 167184
$ws.Range("L66").Value = 71996.9985
$ws.Range("M66").Value = -163440
$ws.Range("N66").Value = -79484.9985
$ws.Range("H99").Value = 4598.2
$ws.Range("I99").Value = 1327.6666
$ws.Range("J99").Value = 5999.857
$ws.Range("K99").Value = 3982.9998
$ws.Range("L99").Value = 17999.571
$ws.Range("M99").Value = -1736.9998
$ws.Range("N99").Value = -22491.571
$ws.Range("H135").Value = 2553
$ws.Range("I135").Value = 1919.7
$ws.Range("J135").Value = 4664
$ws.Range("K135").Value = 17277.3
$ws.Range("L135").Value = 41976
$ws.Range("M135").Value = -14742.3
$ws.Range("N135").Value = -47046

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2149.4614
$ws.Range("I122").Value = 1758.0588
$ws.Range("J122").Value = 2888.7778
$ws.Range("K122").Value = 5274.1764
$ws.Range("L122").Value = 8666.3334
$ws.Range("M122").Value = -2824.1764
$ws.Range("N122").Value = -13566.3334
$ws.Range("H136").Value = 23640.133
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 23640.133
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 70920.399
$ws.Range("N136").Value = -76020.399

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2174.5715
$ws.Range("I22").Value = 854.5
$ws.Range("J22").Value = 3934.6667
$ws.Range("K22").Value = 854.5
$ws.Range("L22").Value = 3934.6667
$ws.Range("M22").Value = -559.5
$ws.Range("N22").Value = -4524.6667
$ws.Range("H27").Value = 2174.5715
$ws.Range("I27").Value = 854.5
$ws.Range("J27").Value = 3934.6667
$ws.Range("K27").Value = 854.5
$ws.Range("L27").Value = 3934.6667
$ws.Range("M27").Value = -747.5
$ws.Range("N27").Value = -4148.6667
$ws.Range("H46").Value = 2360.5557
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2360.5557
$ws.Range("K46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("M46").Value = 2360.5557
$ws.Range("N46").Value = -2736.5557
$ws.Range("H61").Value = 1693.6
$ws.Range("I61").Value = 1885.9375
$ws.Range("J61").Value = 924.25
$ws.Range("K61").Value = 1885.9375
$ws.Range("L61").Value = 924.25
$ws.Range("M61").Value = -1683.9375
$ws.Range("N61").Value = -1328.25
$ws.Range("H113").Value = 1693.6
$ws.Range("I113").Value = 1885.9375
$ws.Range("J113").Value = 924.25
$ws.Range("K113").Value = 1885.9375
$ws.Range("L113").Value = 924.25
$ws.Range("M113").Value = 284.0625
$ws.Range("N113").Value = -5264.25

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 200
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -536
$ws.Range("H62").Value = 8019
$ws.Range("I62").Value = 6996.6665
$ws.Range("J62").Value = 8457.143
$ws.Range("K62").Value = 6996.6665
$ws.Range("L62").Value = 8457.143
$ws.Range("M62").Value = -6372.6665
$ws.Range("N62").Value = -9705.143
$ws.Range("H65").Value = 8019
$ws.Range("I65").Value = 6996.6665
$ws.Range("J65").Value = 8457.143
$ws.Range("K65").Value = 34983.3325
$ws.Range("L65").Value = 42285.715
$ws.Range("M65").Value = -31863.3325
$ws.Range("N65").Value = -48525.715
